$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the group labels (Ambient/Reduced) and their corresponding
# Tukey compact-letter-display labels (a/b) between rows 2 and 3,
# matching the re-run of the Tukey post-hoc test on the new dataset.
$ws.Range("A2").Value = "Reduced"
$ws.Range("B2").Value = "b"
$ws.Range("A3").Value = "Ambient"
$ws.Range("B3").Value = "a"
